# Update the Code Smells metrics sheet: refreshed LOC_class / WMC_class
# figures (LOC_class 50 -> 70, WMC_class 1 -> 13) for the MongoWorker
# methods, plus corrected LOC_method / CYCLO_method for setCloudClient()
# and WMC_class for the inner class row.
#
# The source cells are stored as text (not numbers), so each write
# temporarily switches the cell to text format, assigns the value, then
# restores the cell's original number format to avoid leaving stray
# formatting behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $originalFormat = $cell.NumberFormat
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = $originalFormat
}

# LOC_class 50 -> 70, WMC_class 1 -> 13 for rows 2-8 (MongoWorker methods)
Set-TextValue "F2" "70"
Set-TextValue "G2" "13"

Set-TextValue "F3" "70"
Set-TextValue "G3" "13"
Set-TextValue "I3" "23"
Set-TextValue "J3" "6"

Set-TextValue "F4" "70"
Set-TextValue "G4" "13"

Set-TextValue "F5" "70"
Set-TextValue "G5" "13"

Set-TextValue "F6" "70"
Set-TextValue "G6" "13"

Set-TextValue "F7" "70"
Set-TextValue "G7" "13"

Set-TextValue "F8" "70"
Set-TextValue "G8" "13"

# WMC_class 0 -> 1 for the MongoWorker.InnerClass row
Set-TextValue "G9" "1"
